$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Value = "'26.315.62"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.Value = "'  -4.40%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.Value = "'1.757.24"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.Value = "'  -4.10%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 4)
$c.Value = "'1.002"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.Value = "'  +0.06%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.Value = "'303.23"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.Value = "'  -2.91%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.4279"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.Value = "'  +0.60%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.3602"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.Value = "'  -1.27%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.06987"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.Value = "'  -3.83%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.8270"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.Value = "'  -4.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.Value = "'20.02"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.Value = "'  -2.89%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.Value = "'1.734.18"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.Value = "'  -5.28%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.Value = "'5.197"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.Value = "'  -3.73%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.Value = "'6.338"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.Value = "'  -2.56%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.Value = "'0.06783"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.Value = "'  -2.19%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.Value = "'1.005"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.Value = "'  +0.22%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.Value = "'78.82"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.Value = "'  -2.03%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.Value = "'0.000008633"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.Value = "'  -2.56%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.Value = "'1.003"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.Value = "'  +0.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.Value = "'14.86"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.Value = "'  -3.41%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.Value = "'26.247.31"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.Value = "'  -4.61%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.Value = "'4.959"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.Value = "'  -3.66%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.Value = "'11.02"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.Value = "'  +1.83%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.Value = "'1.955.98"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.Value = "'  -4.83%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.Value = "'1.904"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.Value = "'151.79"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.Value = "'  -1.76%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.Value = "'18.04"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.Value = "'  -4.00%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.Value = "'114.44"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.Value = "'  +0.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.Value = "'4.997"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.Value = "'  -2.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.Value = "'1.659"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.Value = "'  -8.57%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 4)
$c.Value = "'0.08883"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.Value = "'  +0.54%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 4)
$c.Value = "'0.7146"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.Value = "'  -4.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.Value = "'4.284"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.Value = "'  -5.35%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.Value = "'1.088"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.Value = "'  -3.78%  "
$c.Style = "Normal"
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Cells.Item(35, 4)
$c.Value = "'2.758"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.Value = "'  -7.55%  "
$c.Style = "Normal"
$ws.Cells.Item(36, 2).Value = "Frax"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Cells.Item(36, 4)
$c.Value = "'1.001"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.Value = "'  -0.02%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.Value = "'1.060"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.Value = "'  -2.72%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.05069"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.Value = "'  -4.36%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.01873"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.Value = "'  -3.01%  "
$c.Style = "Normal"
$ws.Cells.Item(40, 2).Value = "TheSandbox"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.4869"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.Value = "'  -3.94%  "
$c.Style = "Normal"
$ws.Cells.Item(41, 2).Value = "Algorand"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.1593"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.Value = "'  -3.74%  "
$c.Style = "Normal"
$ws.Cells.Item(42, 2).Value = "MXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Cells.Item(42, 4)
$c.Value = "'2.462"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.Value = "'  -11.97%  "
$c.Style = "Normal"
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Cells.Item(43, 4)
$c.Value = "'6.107"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.Value = "'  -5.77%  "
$c.Style = "Normal"
$ws.Cells.Item(44, 2).Value = "PaxosStandard"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$c = $ws.Cells.Item(44, 4)
$c.Value = "'1.003"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.Value = "'  +0.04%  "
$c.Style = "Normal"
$ws.Cells.Item(45, 2).Value = "Aptos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'7.865"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.Value = "'  -5.33%  "
$c.Style = "Normal"
$ws.Cells.Item(46, 2).Value = "Quant"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Cells.Item(46, 4)
$c.Value = "'104.34"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.Value = "'  -0.86%  "
$c.Style = "Normal"
$ws.Cells.Item(47, 2).Value = "PaxDollar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Cells.Item(47, 4)
$c.Value = "'1.002"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.Value = "'  +0.15%  "
$c.Style = "Normal"
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'9.919"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.Value = "'  -4.78%  "
$c.Style = "Normal"
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(49, 4)
$c.Value = "'0.06156"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.Value = "'  -4.88%  "
$c.Style = "Normal"
$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.4445"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.Value = "'  -4.87%  "
$c.Style = "Normal"
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Cells.Item(51, 4)
$c.Value = "'1.559"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.Value = "'  -3.29%  "
$c.Style = "Normal"
